$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2,1).Value2 = 9999
$excel.CalculateFull()
Write-Host "done"
